$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "Regenerate"
$ws.Range("B9").Value = "N"
$ws.Range("C9").Value = "N/A"
$ws.Range("D9").Value = "治疗一次角色"

$ws.Range("H24").Select()
